$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, pushing the existing rows 214:347 down to 215:348.
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with a new price record (same constant
# columns as the surrounding rows in this sub-dataset), dated 2022-08-04
# (serial 44777) with a volume of 120.
$ws.Range("A214").Value = 3
$ws.Range("B214").Value = 'Femacal de La Calera'
$ws.Range("C214").Value = 'Coquimbo'
$ws.Range("D214").Value = 44777
$ws.Range("E214").Value = 5
$ws.Range("F214").Value = 100112039
$ws.Range("G214").Value = 'Ciboulette'
$ws.Range("H214").Value = 'Sin especificar'
$ws.Range("I214").Value = 'Primera'
$ws.Range("J214").Value = 120
$ws.Range("K214").Value = 1500
$ws.Range("L214").Value = 1500
$ws.Range("M214").Value = 1500
$ws.Range("N214").Value = '$/docena de atados'
$ws.Range("O214").Value = 'Provincia de Quillota'
$ws.Range("P214").Value = 500
$ws.Range("Q214").Value = 3
$ws.Range("R214").Value = 'Hortaliza'
